$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 154 ("「結ばれしター」ة ..." post) entirely.
# This shifts all subsequent rows (155-192) up by one (154-191),
# and reduces the used range from A1:C192 to A1:C191.
$ws.Rows.Item(154).Delete()
